# "Fix, New feature #6"
# - Add a new "ChromeProfile" column (H) to the Instgram sheet.
# - Replace the sample row's Profile_Url / Profile_Name with a new profile.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added column H.
$ws.Cells.Item(1, 8).Value = "ChromeProfile"

# Swap in the new sample profile on row 2.
$ws.Range("A2").Value = "https://www.instagram.com/_heismannu_"
$ws.Range("B2").Value = "Koe"
